$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Select row 4 (entire row) and insert a new blank row there,
# shifting the existing row 4 (TestCheckin0003...) and below down by one.
$row = $ws.Rows.Item(4)
$row.Select()
$row.Insert()

# Match the resulting selection state from the saved file (active cell A4,
# selection spanning the full inserted row).
$ws.Rows.Item(4).Select()
